# edit.ps1 -- LAB-3: Modified the symbol table to generate the positions for symbols
#
# 1) Turns the HYPERLINK field-code in paragraph 1 into a real w:hyperlink run.
# 2) Fixes "collisions" -> "collisions." in paragraph 2 (this also clears the
#    gramStart/gramEnd proofErr markers Word had attached to the old text).
# 3) Appends four new paragraphs describing the hash function / add / get /
#    remove behaviour of the symbol table.

$d = $word.ActiveDocument

# --- 1) Paragraph 1: convert the field-code hyperlink into a real w:hyperlink ---
$field = $d.Fields.Item(1)
$field.Delete()

$hypRange = $d.Paragraphs.Item(1).Range
$d.Hyperlinks.Add($hypRange, "https://github.com/Gabarsolon/FLCD/tree/main/Lab2/src") | Out-Null

# --- 2) Paragraph 2: "collisions" -> "collisions.", proofErr markers go away ---
# Deleting the whole paragraph (incl. its end-of-paragraph mark) drops the
# proofErr annotations Word had recorded for it; we then retype the sentence.
$oldPara2 = $d.Paragraphs.Item(2).Range
$oldPara2.Delete()

$insStart = $d.Paragraphs.Item(2).Range.Start
$insA = $d.Range($insStart, $insStart)
$insA.InsertAfter("The Symbol Table which I’ve implemented uses a hash table with separate chaining for ")

$afterAEnd = $d.Paragraphs.Item(2).Range.End
$insB = $d.Range($afterAEnd - 1, $afterAEnd - 1)

# --- 3) Append the four new paragraphs right after paragraph 2 ---
$insB.InsertAfter("collisions." + [char]13 + "An element is inserted to a bucket using this hash function: h(k) = hashCode(k) % numberOfBuckets, where the numberOfBuckets Is a prime number. When the ratio between the number of elements and the number of buckets is greater than 0.75 (threshold) we changed the capacity to the next prime number and rehash every key." + [char]13 + "The add function puts an entry to its corresponding bucket according to its key if its not already existing and then returns its associated value." + [char]13 + "The get function returns the value associated to a key by finding its bucket and then the node from linked list." + [char]13 + "(optional) The remove function deletes an entry from the hash table and assures that the links between the elements are consistent")

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("Para " + $i + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}
